$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.075.95'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '1.651.59'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '217.39'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '0.5268'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.2594'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").Value = '0.06312'
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").Value = '20.33'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D11").Value = '0.07796'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '4.516'
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").Value = '1.640.72'
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").Value = '1.877.33'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = '0.5478'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '0.0₅8175'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '65.45'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '26.069.85'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '4.574'
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("D21").Value = '190.64'
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").Value = '10.08'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '143.49'
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").Value = '7.217'
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").Value = '15.98'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '0.05809'
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").Value = '1.272'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '3.543'
$ws.Range("D33").Value = '3.262'
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").Value = '1.594'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").Value = '2.796'
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '0.9413'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").Value = '0.5748'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("D39").Value = '0.01604'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '105.37'
$ws.Range("E40").Value = '  +4.88%  '
$ws.Range("D41").Value = '0.8490'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D43").Value = '5.711'
$ws.Range("E43").Value = '  -4.62%  '
$ws.Range("D44").Value = '1.029.19'
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("D45").Value = '1.793.74'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("D46").Value = '57.11'
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").Value = '0.4329'
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("D49").Value = '0.05143'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("D50").Value = '7.829'
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("D51").Value = '1.449'
$ws.Range("E51").Value = '  -0.45%  '
